$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on D:E so numeric-looking strings
# (e.g. "98.00", "2.10") are preserved exactly as typed instead of
# being auto-coerced to numbers by the input parser.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range('D2').Value = '42.979.54'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '2.304.76'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '305.41'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').Value = '98.00'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  -1.62%  '
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').Value = '35.97'
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').Value = '18.24'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').Value = '0.119'
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('D14').Value = '6.80'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').Value = '2.664.37'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '2.306.63'
$ws.Range('E16').Value = '  -2.37%  '
$ws.Range('D17').Value = '0.785'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '42.919.75'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '12.62'
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D22').Value = '68.01'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').Value = '236.54'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  -1.64%  '
$ws.Range('D25').Value = '2.49'
$ws.Range('E25').Value = '  +2.58%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').Value = '25.51'
$ws.Range('E28').Value = '  +3.28%  '
$ws.Range('D29').Value = '165.78'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '33.46'
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '9.08'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '4.83'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E35').Value = '  -2.89%  '
$ws.Range('D36').Value = '17.05'
$ws.Range('E36').Value = '  -5.81%  '
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('E40').Value = '  -1.12%  '
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').Value = '2.009.43'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('E44').Value = '  -1.93%  '
$ws.Range('D45').Value = '10.07'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = '17.84'
$ws.Range('E46').Value = '  +2.46%  '
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('D48').Value = '2.80'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('E49').Value = '  +3.61%  '
$ws.Range('D50').Value = '53.73'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').Value = '2.531.41'
$ws.Range('E51').Value = '  +0.09%  '

# Restore the original (default) cell style now that the text values
# are committed, so formatting matches the source workbook.
$numRange.Style = "Normal"
